# Apply stock report corrections: several quantity (F) and value (G) cells
# were recalculated, which also changes several company Sub Total (B) cells,
# and the KUS-Floor Wiper rows (518/519) had their data swapped.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "F21" = 0
    "G21" = 0
    "B22" = 0
    "F24" = 4
    "G24" = 108.6
    "B46" = 27271.59
    "F48" = 32
    "G48" = 6296.64
    "F55" = 213
    "G55" = 41085.57
    "B85" = 157512.47
    "F137" = 52
    "G137" = 4615.52
    "F143" = 41
    "G143" = 1634.26
    "F152" = 58
    "G152" = 2460.94
    "B159" = 72150.57000000001
    "F174" = 237
    "G174" = 4664.16
    "F178" = 137
    "G178" = 4386.74
    "F179" = 80
    "G179" = 3561.6
    "B180" = 36677.15
    "F188" = 32
    "G188" = 2972.48
    "F191" = 58
    "G191" = 7338.74
    "B198" = 44405.82
    "F259" = 62
    "G259" = 436.48
    "B267" = 21237.1
    "F287" = 19
    "G287" = 4384.44
    "F288" = 130
    "G288" = 3218.8
    "B298" = 109944.05
    "F314" = 84
    "G314" = 1697.64
    "F329" = 204
    "G329" = 6793.2
    "B349" = 141559.43
    "F361" = 54
    "G361" = 3980.34
    "F382" = 28
    "G382" = 3074.96
    "F390" = 12
    "G390" = 1178.4
    "F398" = 1
    "G398" = 154.66
    "F399" = 131
    "G399" = 7689.7
    "F409" = 171
    "G409" = 29297.43
    "F419" = 286
    "G419" = 11777.48
    "F420" = 79
    "G420" = 7317.77
    "F422" = 14
    "G422" = 2011.24
    "B423" = 148687.73
    "F426" = 21
    "G426" = 3855.18
    "F436" = 4
    "G436" = 894.4
    "B437" = 21657.42
    "F445" = 88
    "G445" = 4487.12
    "B455" = 24072.9
    "F463" = 32
    "G463" = 2493.44
    "F476" = 291
    "G476" = 5933.49
    "B481" = 43869.46
    "F487" = 6
    "G487" = 331.98
    "F494" = 99
    "G494" = 5544.99
    "F496" = 99
    "G496" = 14876.73
    "B497" = 38955.83
    "B518" = 58047
    "D518" = 105.54
    "E518" = 126.1
    "F518" = 69
    "G518" = 7282.26
    "B519" = 47097
    "D519" = 112.28
    "E519" = 134.16
    "F519" = 18
    "G519" = 2021.04
    "F560" = 114
    "G560" = 5762.7
    "F564" = 177
    "G564" = 5000.25
    "B567" = 51640.8
    "F591" = 314
    "G591" = 11608.58
    "F592" = 90
    "G592" = 2817
    "F595" = 73
    "G595" = 3655.84
    "F596" = 175
    "G596" = 5078.5
    "F597" = 126
    "G597" = 7767.9
    "B610" = 56812.24
    "F621" = 242
    "G621" = 14694.24
    "F625" = 29
    "G625" = 1865.28
    "B638" = 147483.87
    "F663" = 66
    "G663" = 5200.8
    "B667" = 26295.24
    "F672" = 71
    "G672" = 18892.39
    "F674" = 97
    "G674" = 5061.46
    "F678" = 17
    "G678" = 798.3200000000001
    "B688" = 87748.37
    "F715" = 47
    "G715" = 6135.85
    "F718" = 129
    "G718" = 3508.8
    "F719" = 107
    "G719" = 2910.4
    "B720" = 30845.39
    "F745" = 93
    "G745" = 5756.7
    "F758" = 11
    "G758" = 9118.120000000001
    "B773" = 135174.2
    "F778" = 11
    "G778" = 787.16
    "B785" = 15325.54
    "F810" = 71
    "G810" = 6466.68
    "B815" = 38415.53
    "F822" = 46
    "G822" = 5088.98
    "F825" = 22
    "G825" = 820.16
    "F830" = 260
    "G830" = 39293.8
    "F831" = 28
    "G831" = 7656.88
    "B837" = 194294.6
    "F839" = 31
    "G839" = 7733.88
    "F843" = 69
    "G843" = 7507.89
    "F849" = 33
    "G849" = 4828.89
    "F862" = 15
    "G862" = 708.15
    "B867" = 203432.96
    "F892" = 2
    "G892" = 1147.12
    "B904" = 38470.05
    "B923" = 2647290.47
    "B924" = 2647290.47
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
